# Update TPM (transcripts per million) derived NATMI ligand-receptor
# statistics for the Vcam1-Itgad pair after re-running the analysis
# scripts with new TPM values. Also fixes the Target cluster label for
# the Inflammatory-Mac rows, which should read Neutrophils.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value2 = "Neutrophils"
$ws.Range("G2").Value2 = 9.6736855
$ws.Range("H2").Value2 = 19.347371
$ws.Range("I2").Value2 = 0.115977165622779
$ws.Range("J2").Value2 = 0.09051257664205575
$ws.Range("K2").Value2 = 3
$ws.Range("L2").Value2 = 1
$ws.Range("M2").Value2 = 0.383214
$ws.Range("N2").Value2 = 1.149642
$ws.Range("O2").Value2 = 0.891270657753252
$ws.Range("P2").Value2 = 0.891270657753252
$ws.Range("Q2").Value2 = 3.707091715197
$ws.Range("R2").Value2 = 22.242550291182
$ws.Range("S2").Value2 = 0.1033670446889721
$ws.Range("T2").Value2 = 0.08067120371870666

# Row 3
$ws.Range("G3").Value2 = 9.6736855
$ws.Range("H3").Value2 = 19.347371
$ws.Range("I3").Value2 = 0.115977165622779
$ws.Range("J3").Value2 = 0.09051257664205575
$ws.Range("K3").Value2 = 2
$ws.Range("L3").Value2 = 0.6666666666666666
$ws.Range("M3").Value2 = 0.04674966666666667
$ws.Range("N3").Value2 = 0.140249
$ws.Range("O3").Value2 = 0.108729342246748
$ws.Range("P3").Value2 = 0.108729342246748
$ws.Range("Q3").Value2 = 0.4522415725631667
$ws.Range("R3").Value2 = 2.713449435379
$ws.Range("S3").Value2 = 0.01261012093380691
$ws.Range("T3").Value2 = 0.009841372923349087

# Row 4
$ws.Range("D4").Value2 = "Neutrophils"
$ws.Range("G4").Value2 = 40.41312266666666
$ws.Range("I4").Value2 = 0.484510212870336
$ws.Range("J4").Value2 = 0.5671926996248948
$ws.Range("K4").Value2 = 3
$ws.Range("L4").Value2 = 1
$ws.Range("M4").Value2 = 0.383214
$ws.Range("N4").Value2 = 1.149642
$ws.Range("O4").Value2 = 0.891270657753252
$ws.Range("P4").Value2 = 0.891270657753252
$ws.Range("Q4").Value2 = 15.486874389584
$ws.Range("R4").Value2 = 139.381869506256
$ws.Range("S4").Value2 = 0.4318297361131125
$ws.Range("T4").Value2 = 0.5055222104675227

# Row 5
$ws.Range("G5").Value2 = 40.41312266666666
$ws.Range("I5").Value2 = 0.484510212870336
$ws.Range("J5").Value2 = 0.5671926996248948
$ws.Range("K5").Value2 = 2
$ws.Range("L5").Value2 = 0.6666666666666666
$ws.Range("M5").Value2 = 0.04674966666666667
$ws.Range("N5").Value2 = 0.140249
$ws.Range("O5").Value2 = 0.108729342246748
$ws.Range("P5").Value2 = 0.108729342246748
$ws.Range("Q5").Value2 = 1.889300013625778
$ws.Range("R5").Value2 = 17.003700122632
$ws.Range("S5").Value2 = 0.05268047675722348
$ws.Range("T5").Value2 = 0.06167048915737211

# Row 6
$ws.Range("D6").Value2 = "Neutrophils"
$ws.Range("G6").Value2 = 1.927632333333333
$ws.Range("H6").Value2 = 5.782896999999999
$ws.Range("I6").Value2 = 0.0231102545542569
$ws.Range("J6").Value2 = 0.02705405855532591
$ws.Range("K6").Value2 = 3
$ws.Range("L6").Value2 = 1
$ws.Range("M6").Value2 = 0.383214
$ws.Range("N6").Value2 = 1.149642
$ws.Range("O6").Value2 = 0.891270657753252
$ws.Range("P6").Value2 = 0.891270657753252
$ws.Range("Q6").Value2 = 0.7386956969859999
$ws.Range("R6").Value2 = 6.648261272874
$ws.Range("S6").Value2 = 0.02059749177741763
$ws.Range("T6").Value2 = 0.02411248856350031

# Row 7
$ws.Range("G7").Value2 = 1.927632333333333
$ws.Range("H7").Value2 = 5.782896999999999
$ws.Range("I7").Value2 = 0.0231102545542569
$ws.Range("J7").Value2 = 0.02705405855532591
$ws.Range("K7").Value2 = 2
$ws.Range("L7").Value2 = 0.6666666666666666
$ws.Range("M7").Value2 = 0.04674966666666667
$ws.Range("N7").Value2 = 0.140249
$ws.Range("O7").Value2 = 0.108729342246748
$ws.Range("P7").Value2 = 0.108729342246748
$ws.Range("Q7").Value2 = 0.09011616903922222
$ws.Range("R7").Value2 = 0.811045521353
$ws.Range("S7").Value2 = 0.002512762776839264
$ws.Range("T7").Value2 = 0.002941569991825591

# Row 8
$ws.Range("D8").Value2 = "Neutrophils"
$ws.Range("G8").Value2 = 26.8036935
$ws.Range("H8").Value2 = 53.607387
$ws.Range("I8").Value2 = 0.3213476808142776
$ws.Range("J8").Value2 = 0.2507908037954016
$ws.Range("K8").Value2 = 3
$ws.Range("L8").Value2 = 1
$ws.Range("M8").Value2 = 0.383214
$ws.Range("N8").Value2 = 1.149642
$ws.Range("O8").Value2 = 0.891270657753252
$ws.Range("P8").Value2 = 0.891270657753252
$ws.Range("Q8").Value2 = 10.271550600909
$ws.Range("R8").Value2 = 61.62930360545401
$ws.Range("S8").Value2 = 0.2864077588468232
$ws.Range("T8").Value2 = 0.2235224846571944

# Row 9
$ws.Range("G9").Value2 = 26.8036935
$ws.Range("H9").Value2 = 53.607387
$ws.Range("I9").Value2 = 0.3213476808142776
$ws.Range("J9").Value2 = 0.2507908037954016
$ws.Range("K9").Value2 = 2
$ws.Range("L9").Value2 = 0.6666666666666666
$ws.Range("M9").Value2 = 0.04674966666666667
$ws.Range("N9").Value2 = 0.140249
$ws.Range("O9").Value2 = 0.108729342246748
$ws.Range("P9").Value2 = 0.108729342246748
$ws.Range("Q9").Value2 = 1.2530637365605
$ws.Range("R9").Value2 = 7.518382419363001
$ws.Range("S9").Value2 = 0.03493992196745431
$ws.Range("T9").Value2 = 0.02726831913820725

# Row 10
$ws.Range("D10").Value2 = "Neutrophils"
$ws.Range("G10").Value2 = 2.854607333333333
$ws.Range("H10").Value2 = 8.563822
$ws.Range("I10").Value2 = 0.03422369555905032
$ws.Range("J10").Value2 = 0.04006402705173345
$ws.Range("K10").Value2 = 3
$ws.Range("L10").Value2 = 1
$ws.Range("M10").Value2 = 0.383214
$ws.Range("N10").Value2 = 1.149642
$ws.Range("O10").Value2 = 0.891270657753252
$ws.Range("P10").Value2 = 0.891270657753252
$ws.Range("Q10").Value2 = 1.093925494636
$ws.Range("R10").Value2 = 9.845329451724
$ws.Range("S10").Value2 = 0.03050257565166183
$ws.Range("T10").Value2 = 0.03570789174264256

# Row 11
$ws.Range("G11").Value2 = 2.854607333333333
$ws.Range("H11").Value2 = 8.563822
$ws.Range("I11").Value2 = 0.03422369555905032
$ws.Range("J11").Value2 = 0.04006402705173345
$ws.Range("K11").Value2 = 2
$ws.Range("L11").Value2 = 0.6666666666666666
$ws.Range("M11").Value2 = 0.04674966666666667
$ws.Range("N11").Value2 = 0.140249
$ws.Range("O11").Value2 = 0.108729342246748
$ws.Range("P11").Value2 = 0.108729342246748
$ws.Range("Q11").Value2 = 0.1334519412975556
$ws.Range("R11").Value2 = 1.201067471678
$ws.Range("S11").Value2 = 0.003721119907388491
$ws.Range("T11").Value2 = 0.004356135309090896

# Row 12
$ws.Range("D12").Value2 = "Neutrophils"
$ws.Range("G12").Value2 = 1.737518333333333
$ws.Range("H12").Value2 = 5.212555
$ws.Range("I12").Value2 = 0.0208309905793004
$ws.Range("J12").Value2 = 0.02438583433058843
$ws.Range("K12").Value2 = 3
$ws.Range("L12").Value2 = 1
$ws.Range("M12").Value2 = 0.383214
$ws.Range("N12").Value2 = 1.149642
$ws.Range("O12").Value2 = 0.891270657753252
$ws.Range("P12").Value2 = 0.891270657753252
$ws.Range("Q12").Value2 = 0.6658413505899999
$ws.Range("R12").Value2 = 5.99257215531
$ws.Range("S12").Value2 = 0.01856605067526487
$ws.Range("T12").Value2 = 0.02173437860368539

# Row 13
$ws.Range("G13").Value2 = 1.737518333333333
$ws.Range("H13").Value2 = 5.212555
$ws.Range("I13").Value2 = 0.0208309905793004
$ws.Range("J13").Value2 = 0.02438583433058843
$ws.Range("K13").Value2 = 2
$ws.Range("L13").Value2 = 0.6666666666666666
$ws.Range("M13").Value2 = 0.04674966666666667
$ws.Range("N13").Value2 = 0.140249
$ws.Range("O13").Value2 = 0.108729342246748
$ws.Range("P13").Value2 = 0.108729342246748
$ws.Range("Q13").Value2 = 0.08122840291055555
$ws.Range("R13").Value2 = 0.7310556261950001
$ws.Range("S13").Value2 = 0.002264939904035537
$ws.Range("T13").Value2 = 0.002651455726903046

